$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldName = $ws.Name
$newName = "Leerling Iwan"

# Rename the worksheet
$ws.Name = $newName

# Reset the print area so the stored defined name references the new sheet name
$ws.PageSetup.PrintArea = '$AD$1:$AK$11'

# Update all chart series formulas that still reference the old sheet name
for ($i = 1; $i -le $ws.ChartObjects().Count; $i++) {
    $co = $ws.ChartObjects($i)
    $chart = $co.Chart
    for ($j = 1; $j -le $chart.SeriesCollection().Count; $j++) {
        $s = $chart.SeriesCollection($j)
        $f = $s.Formula
        if ($f -like "*$oldName*") {
            if ($f -match "SERIES\(,,'(\(.*\)),1\)") {
                # Special-case series using a parenthesised multi-range union
                # reference (no name/category part). Reading/writing .Values
                # directly with the un-rewritten getter value mangles the
                # stored formula, so rebuild the range expression from the
                # .Formula text instead and assign it via .Values.
                $rangeExpr = $Matches[1]
                $newRangeExpr = $rangeExpr.Replace($oldName, $newName)
                $s.Values = $newRangeExpr
            }
            else {
                $newF = $f.Replace($oldName, $newName)
                $s.Formula = $newF
            }
        }
    }
}
